$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Folder Name"
$ws.Range("B1").Value = "Timestamp"

# Update data rows 2 and 3
$ws.Range("A2").Value = "vishal kotnod_181"
$ws.Range("B2").Value = "2025-04-25 08:01:59"

$ws.Range("A3").Value = "vishal kotnod_181"
$ws.Range("B3").Value = "2025-04-26 15:28:38"

# Remove now-unused columns C and D, and rows 4-5
$ws.Range("C1:D5").Clear()
$ws.Range("A4:B5").Clear()
